# The "kernel_summary" sheet contains, per thread-count variant, a set of
# columns: Count, Mean Kernel Time, Percentage, (Cumulative Percentage for
# 256thread), Total Kernel Time, and Category. This edit removes the
# "Category" column for both the 256thread (column M) and 512thread
# (column R) variants, simplifying the rocprof comparison data.
#
# Deleting column M shifts 512thread::Count (N) -> M, 512thread::Mean
# Kernel Time (O) -> N, 512thread::Percentage (P) -> O and
# 512thread::Total Kernel Time (Q) -> P, while the former 512thread::Category
# column (R) ends up as Q and is then removed as well, leaving the sheet
# with columns A:P.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("kernel_summary")

# Delete the right-most "Category" column first so column letters for the
# left-most deletion remain valid.
$ws.Range("R1").EntireColumn.Delete()
$ws.Range("M1").EntireColumn.Delete()
